$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (round.vhd): update the MATLAB-implementation formula text in D5
#   (@rca.m) +  3*(@muxnto1_nbit.m) + (@orn)  ->  ... + (@orn.m)
$ws.Range("D5").Value = "(@rca.m) +  3*(@muxnto1_nbit.m) + (@orn.m)"

# Row 6 (relu.vhd): fill in Description (C6), MATLAB implementation (D6), Done? (F6)
$c6 = $ws.Range("C6")
$c6.Value = "Combinational logic capable of performing a quantization,`n the input and output parallelism is generic."
$c6.NumberFormat = "@"
$c6.HorizontalAlignment = -4108
$c6.VerticalAlignment = -4108
$c6.WrapText = $true

$d6 = $ws.Range("D6")
$d6.Value = "(@orn.m) + 3*(@nand.m) + (@nor.m) + n_O*(@muxnto1_nbit.m) "
$d6.HorizontalAlignment = -4108
$d6.VerticalAlignment = -4108

$f6 = $ws.Range("F6")
$f6.Value = "✔"
$f6.HorizontalAlignment = -4108
$f6.VerticalAlignment = -4108

# Row 8 (pool.vhd): fill in Description (C8), MATLAB implementation (D8), Done? (F8)
$c8 = $ws.Range("C8")
$c8.Value = "Sequential circuit that performs the sum of four quantities in`nthree clock cycle, the data width is generic."
$c8.NumberFormat = "@"
$c8.HorizontalAlignment = -4108
$c8.VerticalAlignment = -4108
$c8.WrapText = $true

$d8 = $ws.Range("D8")
$d8.Value = "(@rca.m)(N+1) + (@rca.m)(N+2) + 2*(@register.m)(N+1) +`n+ (@register.m)(N+1) + 3*(@muxnto1_nbit.m)"
$d8.HorizontalAlignment = -4108
$d8.VerticalAlignment = -4108
$d8.WrapText = $true

$f8 = $ws.Range("F8")
$f8.Value = "✔"
$f8.HorizontalAlignment = -4108
$f8.VerticalAlignment = -4108

# Row 13: add an (empty) formatted cell at C13 (underline font, centered, text-number-format)
$c13 = $ws.Range("C13")
$c13.NumberFormat = "@"
$c13.HorizontalAlignment = -4108
$c13.VerticalAlignment = -4108
$c13.Font.Underline = 2

